$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp in title (row 1)
$ws.Cells.Item(1,1).Value = "Datos actualizados a 20 de Junio de 2020 a las 12:16"

# Row 26: Belgica
$ws.Cells.Item(26,1).Value = "Belgica"
$ws.Cells.Item(26,2).Value = 60550
$ws.Cells.Item(26,3).Value = 74
$ws.Cells.Item(26,4).Value = 16771
$ws.Cells.Item(26,5).Value = 34083
$ws.Cells.Item(26,6).Value = 0
$ws.Cells.Item(26,7).Value = 1
$ws.Cells.Item(26,8).Value = 9696

# Row 41: Oman
$ws.Cells.Item(41,1).Value = "Oman"
$ws.Cells.Item(41,2).Value = 28566
$ws.Cells.Item(41,3).Value = 896
$ws.Cells.Item(41,4).Value = 14780
$ws.Cells.Item(41,5).Value = 13658
$ws.Cells.Item(41,6).Value = 0
$ws.Cells.Item(41,7).Value = 3
$ws.Cells.Item(41,8).Value = 128

# Row 42: Filipinas
$ws.Cells.Item(42,1).Value = "Filipinas"
$ws.Cells.Item(42,2).Value = 28459
$ws.Cells.Item(42,3).Value = 0
$ws.Cells.Item(42,4).Value = 7378
$ws.Cells.Item(42,5).Value = 19951
$ws.Cells.Item(42,6).Value = 0
$ws.Cells.Item(42,7).Value = 0
$ws.Cells.Item(42,8).Value = 1130

# Row 43: Afganistan
$ws.Cells.Item(43,1).Value = "Afganistan"
$ws.Cells.Item(43,2).Value = 28424
$ws.Cells.Item(43,3).Value = 546
$ws.Cells.Item(43,4).Value = 8292
$ws.Cells.Item(43,5).Value = 19563
$ws.Cells.Item(43,6).Value = 0
$ws.Cells.Item(43,7).Value = 21
$ws.Cells.Item(43,8).Value = 569

# Row 48: Rumania
$ws.Cells.Item(48,1).Value = "Rumania"
$ws.Cells.Item(48,2).Value = 23730
$ws.Cells.Item(48,3).Value = 330
$ws.Cells.Item(48,4).Value = 16735
$ws.Cells.Item(48,5).Value = 5495
$ws.Cells.Item(48,6).Value = 0
$ws.Cells.Item(48,7).Value = 16
$ws.Cells.Item(48,8).Value = 1500

# Row 50: Barein
$ws.Cells.Item(50,1).Value = "Barein"
$ws.Cells.Item(50,2).Value = 20916
$ws.Cells.Item(50,3).Value = 0
$ws.Cells.Item(50,4).Value = 15287
$ws.Cells.Item(50,5).Value = 5571
$ws.Cells.Item(50,6).Value = 0
$ws.Cells.Item(50,7).Value = 1
$ws.Cells.Item(50,8).Value = 58

# Row 55: Austria
$ws.Cells.Item(55,1).Value = "Austria"
$ws.Cells.Item(55,2).Value = 17323
$ws.Cells.Item(55,3).Value = 52
$ws.Cells.Item(55,4).Value = 16175
$ws.Cells.Item(55,5).Value = 460
$ws.Cells.Item(55,6).Value = 0
$ws.Cells.Item(55,7).Value = 0
$ws.Cells.Item(55,8).Value = 688

# Row 68: Marruecos
$ws.Cells.Item(68,1).Value = "Marruecos"
$ws.Cells.Item(68,2).Value = 9801
$ws.Cells.Item(68,3).Value = 188
$ws.Cells.Item(68,4).Value = 8133
$ws.Cells.Item(68,5).Value = 1455
$ws.Cells.Item(68,6).Value = 0
$ws.Cells.Item(68,7).Value = 0
$ws.Cells.Item(68,8).Value = 213

# Row 70: Malasia
$ws.Cells.Item(70,1).Value = "Malasia"
$ws.Cells.Item(70,2).Value = 8556
$ws.Cells.Item(70,3).Value = 21
$ws.Cells.Item(70,4).Value = 8146
$ws.Cells.Item(70,5).Value = 289
$ws.Cells.Item(70,6).Value = 0
$ws.Cells.Item(70,7).Value = 0
$ws.Cells.Item(70,8).Value = 121

# Row 74: Finlandia
$ws.Cells.Item(74,1).Value = "Finlandia"
$ws.Cells.Item(74,2).Value = 7142
$ws.Cells.Item(74,3).Value = 9
$ws.Cells.Item(74,4).Value = 6200
$ws.Cells.Item(74,5).Value = 616
$ws.Cells.Item(74,6).Value = 0
$ws.Cells.Item(74,7).Value = 0
$ws.Cells.Item(74,8).Value = 326

# Row 77: Consejo Danes para los Refugiados
$ws.Cells.Item(77,1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(77,2).Value = 5672
$ws.Cells.Item(77,3).Value = 195
$ws.Cells.Item(77,4).Value = 807
$ws.Cells.Item(77,5).Value = 4740
$ws.Cells.Item(77,6).Value = 0
$ws.Cells.Item(77,7).Value = 3
$ws.Cells.Item(77,8).Value = 125

# Row 78: Senegal
$ws.Cells.Item(78,1).Value = "Senegal"
$ws.Cells.Item(78,2).Value = 5639
$ws.Cells.Item(78,3).Value = 0
$ws.Cells.Item(78,4).Value = 3788
$ws.Cells.Item(78,5).Value = 1772
$ws.Cells.Item(78,6).Value = 0
$ws.Cells.Item(78,7).Value = 0
$ws.Cells.Item(78,8).Value = 79

# Row 108: Albania
$ws.Cells.Item(108,1).Value = "Albania"
$ws.Cells.Item(108,2).Value = 1891
$ws.Cells.Item(108,3).Value = 53
$ws.Cells.Item(108,4).Value = 1126
$ws.Cells.Item(108,5).Value = 723
$ws.Cells.Item(108,6).Value = 0
$ws.Cells.Item(108,7).Value = 0
$ws.Cells.Item(108,8).Value = 42

# Row 109: Sudan del Sur
$ws.Cells.Item(109,1).Value = "Sudan del Sur"
$ws.Cells.Item(109,2).Value = 1864
$ws.Cells.Item(109,3).Value = 0
$ws.Cells.Item(109,4).Value = 122
$ws.Cells.Item(109,5).Value = 1708
$ws.Cells.Item(109,6).Value = 0
$ws.Cells.Item(109,7).Value = 0
$ws.Cells.Item(109,8).Value = 34

# Row 116: Eslovenia
$ws.Cells.Item(116,1).Value = "Eslovenia"
$ws.Cells.Item(116,2).Value = 1519
$ws.Cells.Item(116,3).Value = 6
$ws.Cells.Item(116,4).Value = 1359
$ws.Cells.Item(116,5).Value = 51
$ws.Cells.Item(116,6).Value = 0
$ws.Cells.Item(116,7).Value = 0
$ws.Cells.Item(116,8).Value = 109

# Row 123: Tunez
$ws.Cells.Item(123,1).Value = "Tunez"
$ws.Cells.Item(123,2).Value = 1156
$ws.Cells.Item(123,3).Value = 10
$ws.Cells.Item(123,4).Value = 1017
$ws.Cells.Item(123,5).Value = 89
$ws.Cells.Item(123,6).Value = 0
$ws.Cells.Item(123,7).Value = 0
$ws.Cells.Item(123,8).Value = 50

# Row 139: Estado de Palestina
$ws.Cells.Item(139,1).Value = "Estado de Palestina"
$ws.Cells.Item(139,2).Value = 708
$ws.Cells.Item(139,3).Value = 33
$ws.Cells.Item(139,4).Value = 437
$ws.Cells.Item(139,5).Value = 268
$ws.Cells.Item(139,6).Value = 0
$ws.Cells.Item(139,7).Value = 0
$ws.Cells.Item(139,8).Value = 3

# Row 187: Namibia
$ws.Cells.Item(187,1).Value = "Namibia"
$ws.Cells.Item(187,2).Value = 46
$ws.Cells.Item(187,3).Value = 1
$ws.Cells.Item(187,4).Value = 19
$ws.Cells.Item(187,5).Value = 27
$ws.Cells.Item(187,6).Value = 0
$ws.Cells.Item(187,7).Value = 0
$ws.Cells.Item(187,8).Value = 0

# Row 202: Fiyi
$ws.Cells.Item(202,1).Value = "Fiyi"
$ws.Cells.Item(202,2).Value = 18
$ws.Cells.Item(202,3).Value = 0
$ws.Cells.Item(202,4).Value = 18
$ws.Cells.Item(202,5).Value = 0
$ws.Cells.Item(202,6).Value = 0
$ws.Cells.Item(202,7).Value = 0
$ws.Cells.Item(202,8).Value = 0

# Row 203: Dominica
$ws.Cells.Item(203,1).Value = "Dominica"
$ws.Cells.Item(203,2).Value = 18
$ws.Cells.Item(203,3).Value = 0
$ws.Cells.Item(203,4).Value = 18
$ws.Cells.Item(203,5).Value = 0
$ws.Cells.Item(203,6).Value = 0
$ws.Cells.Item(203,7).Value = 0
$ws.Cells.Item(203,8).Value = 0

# Row 208: Santa Sede
$ws.Cells.Item(208,1).Value = "Santa Sede"
$ws.Cells.Item(208,2).Value = 12
$ws.Cells.Item(208,3).Value = 0
$ws.Cells.Item(208,4).Value = 12
$ws.Cells.Item(208,5).Value = 0
$ws.Cells.Item(208,6).Value = 0
$ws.Cells.Item(208,7).Value = 0
$ws.Cells.Item(208,8).Value = 0

# Row 209: Islas Turcas y Caicos
$ws.Cells.Item(209,1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(209,2).Value = 12
$ws.Cells.Item(209,3).Value = 0
$ws.Cells.Item(209,4).Value = 11
$ws.Cells.Item(209,5).Value = 0
$ws.Cells.Item(209,6).Value = 0
$ws.Cells.Item(209,7).Value = 0
$ws.Cells.Item(209,8).Value = 1
